# Add a new "11-sep" date column (BN) to the Dataframe Fam sheet, with its
# per-brand counts, matching the existing "BM" column's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell (column BN = 66, row 1) - text, inherits row/header style.
$ws.Cells.Item(1, 66).Value = "11-sep"

# New data values for BN2:BN11 - numeric, centered, integer format (same as
# the existing BM3:BM11 cells).
$values = @(15, 14, 12, 11, 12, 16, 22, 12, 12, 9)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 66)
    $cell.HorizontalAlignment = -4108
    $cell.NumberFormat = "0"
    $cell.Value = $values[$i]
}

# Move the active selection to just below the new column's data, as in the
# saved workbook.
$ws.Range("BN12").Select()
